$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target range stays text-formatted so values are not
# auto-converted to numbers/percentages by Excel when assigned.
$ws.Range("D2:E48").NumberFormat = "@"

$ws.Range("D2").Value = "279.31"
$ws.Range("E2").Value = "6.44%"
$ws.Range("D3").Value = "27.01"
$ws.Range("E3").Value = "1.01%"
$ws.Range("D4").Value = "4.902"
$ws.Range("E4").Value = "4.51%"
$ws.Range("D5").Value = "0.06338"
$ws.Range("E5").Value = "3.89%"
$ws.Range("D6").Value = "6.943"
$ws.Range("E6").Value = "3.60%"
$ws.Range("D7").Value = "3.352"
$ws.Range("E7").Value = "5.94%"
$ws.Range("D8").Value = "0.8829"
$ws.Range("E8").Value = "3.86%"
$ws.Range("D9").Value = "0.9451"
$ws.Range("E9").Value = "4.00%"
$ws.Range("E10").Value = "4.41%"
$ws.Range("D11").Value = "0.05174"
$ws.Range("E11").Value = "2.65%"
$ws.Range("D12").Value = "0.07422"
$ws.Range("E12").Value = "4.46%"
$ws.Range("D13").Value = "0.03142"
$ws.Range("E13").Value = "0.73%"
$ws.Range("D14").Value = "0.09058"
$ws.Range("E14").Value = "0.04%"
$ws.Range("D15").Value = "0.001555"
$ws.Range("E15").Value = "1.10%"
$ws.Range("D16").Value = "0.0006259"
$ws.Range("E16").Value = "1.24%"
$ws.Range("D17").Value = "0.005799"
$ws.Range("E17").Value = "-3.15%"
$ws.Range("D18").Value = "3.476"
$ws.Range("E18").Value = "0.79%"
$ws.Range("D19").Value = "2.279"
$ws.Range("E19").Value = "5.13%"
$ws.Range("D21").Value = "0.1338"
$ws.Range("E21").Value = "2.96%"
$ws.Range("D22").Value = "3.891"
$ws.Range("E22").Value = "-4.88%"
$ws.Range("D23").Value = "0.04311"
$ws.Range("E23").Value = "2.18%"
$ws.Range("D24").Value = "0.001180"
$ws.Range("E24").Value = "0.27%"
$ws.Range("D25").Value = "0.003616"
$ws.Range("E25").Value = "-10.92%"
$ws.Range("E26").Value = "-0.06%"
$ws.Range("D27").Value = "0.0001693"
$ws.Range("E27").Value = "-12.65%"
$ws.Range("D40").Value = "0.04046"
$ws.Range("E40").Value = "2.64%"
$ws.Range("D41").Value = "0.006623"
$ws.Range("E41").Value = "58.09%"
$ws.Range("D42").Value = "0.1166"
$ws.Range("E42").Value = "4.87%"
$ws.Range("D43").Value = "0.002339"
$ws.Range("E43").Value = "10.87%"
$ws.Range("D44").Value = "0.01243"
$ws.Range("E44").Value = "7.56%"
$ws.Range("D45").Value = "0.00005211"
$ws.Range("E45").Value = "2.29%"
$ws.Range("E47").Value = "819.82%"
$ws.Range("D48").Value = "0.02248"
$ws.Range("E48").Value = "5.96%"

# Restore default (Normal) style so no stray number-format style
# is left attached to the edited cells.
$ws.Range("D2:E48").Style = "Normal"

